# Weekly data refresh: insert 3 new rows of data (week of 2023-06-02,
# serial 45079) at the top of the "Ciruela" price table, pushing the
# existing rows (previously 88-168) down to 91-171.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the three new records just above the current row 88.
$ws.Rows("88:90").Insert()

# Shared values for the three new rows.
$mercadoId = 7
$mercado   = "Terminal Hortofrutícola Agro Chillán"
$region    = "Ñuble"
$fecha     = 45079
$codreg    = 16
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria   = "Ciruela"
$unidad    = "$/bandeja 18 kilos granel"
$origen    = "Región de O'Higgins"
$kgUnidad  = 18

# Row 88: Angeleno - Especial
$r = 88
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Angeleno"
$ws.Cells.Item($r, 12).Value = "Especial"
$ws.Cells.Item($r, 13).Value = 40
$ws.Cells.Item($r, 14).Value = 10000
$ws.Cells.Item($r, 15).Value = 10000
$ws.Cells.Item($r, 16).Value = 10000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 556
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 89: Angeleno - Primera
$r = 89
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Angeleno"
$ws.Cells.Item($r, 12).Value = "Primera"
$ws.Cells.Item($r, 13).Value = 30
$ws.Cells.Item($r, 14).Value = 8000
$ws.Cells.Item($r, 15).Value = 8000
$ws.Cells.Item($r, 16).Value = 8000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 444
$ws.Cells.Item($r, 20).Value = $kgUnidad

# Row 90: Angeleno - Segunda
$r = 90
$ws.Cells.Item($r, 1).Value  = $mercadoId
$ws.Cells.Item($r, 2).Value  = $mercado
$ws.Cells.Item($r, 3).Value  = $region
$ws.Cells.Item($r, 4).Value  = $fecha
$ws.Cells.Item($r, 5).Value  = $codreg
$ws.Cells.Item($r, 6).Value  = $tipo
$ws.Cells.Item($r, 7).Value  = $productoId
$ws.Cells.Item($r, 8).Value  = $producto
$ws.Cells.Item($r, 9).Value  = $categoriaId
$ws.Cells.Item($r, 10).Value = $categoria
$ws.Cells.Item($r, 11).Value = "Angeleno"
$ws.Cells.Item($r, 12).Value = "Segunda"
$ws.Cells.Item($r, 13).Value = 30
$ws.Cells.Item($r, 14).Value = 6000
$ws.Cells.Item($r, 15).Value = 6000
$ws.Cells.Item($r, 16).Value = 6000
$ws.Cells.Item($r, 17).Value = $unidad
$ws.Cells.Item($r, 18).Value = $origen
$ws.Cells.Item($r, 19).Value = 333
$ws.Cells.Item($r, 20).Value = $kgUnidad
